$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.676.49"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.752.73"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.51%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.39"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -4.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4341"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -5.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3656"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.36"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07468"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.120"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -2.34%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.65"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -3.25%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.160"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.256"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.749.38"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001068"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "88.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +8.88%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06211"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -7.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.155"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -3.73%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.5306"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.708.18"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.64"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.84%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.323"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -3.95%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.56"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.24%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "153.82"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.357"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.949.27"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.94%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "127.80"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.729"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.50%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09145"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.636"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -9.65%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.67"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02311"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2163"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -5.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.098"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6471"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.79%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06101"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -3.58%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.195"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.421"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.964"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "13.83"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.750"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -3.02%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5938"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -2.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "125.95"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -3.49%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.974"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.67%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06895"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.55%  "
